# Update the "Metadata" worksheet (sheet1) of the StructureDefinition workbook:
#  - Refresh the "Date" property value
#  - Change the "Context" value from "element:Element" to "element:List"
#  - Add a second "Context" row with value "element:Consent"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date property (row 8, column B)
$ws.Range("B8").Value = "2024-03-11T22:11:27+00:00"

# The extension now applies to two contexts instead of one: "List" and "Consent".
# Update the existing Context row to the first new value...
$ws.Range("B20").Value = "element:List"

# ...and duplicate the Context row (to inherit the same style) for the second value.
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Context"
$ws.Range("B21").Value = "element:Consent"

$excel.CutCopyMode = 0
